# Append the latest mod-count entry as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (the "Date" column) so the new
# record lands directly below the existing data, regardless of how many
# rows are already present.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
$newRow = $lastRow + 1

# Write the new values. Assign the date as text first (quote-prefixed)
# so Excel doesn't silently reinterpret "2026/02/16" as a date serial;
# the formatting pass below then restores the normal (non-text) style
# used by the rest of the table, matching the previous row exactly.
$ws.Cells.Item($newRow, 1).Value = "'2026/02/16"
$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1210

# Copy the formatting from the row above so the new row's style matches
# the rest of the table (centered alignment, same number format, etc.)
$ws.Range("A$($lastRow):C$($lastRow)").Copy() | Out-Null
$ws.Range("A$($newRow):C$($newRow)").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
